$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 340. This shifts the existing rows 340:461
# down to 341:462 (dimension grows from A1:R461 to A1:R462), matching the
# target diff which keeps all "old row N" data at "new row N+1" and adds a
# brand-new record at row 340.
$ws.Rows("340").Insert()

# Populate the new row 340 with a new price-report record for Cilantro.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R mirror the record that used to
# sit at row 340 (now shifted to row 341); only D, J, K, L, M, P differ.
$ws.Range("A340").Value = 3
$ws.Range("B340").Value = "Femacal de La Calera"
$ws.Range("C340").Value = "Coquimbo"
$ws.Range("D340").Value = 44900
$ws.Range("E340").Value = 5
$ws.Range("F340").Value = 100112040
$ws.Range("G340").Value = "Cilantro"
$ws.Range("H340").Value = "Sin especificar"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 230
$ws.Range("K340").Value = 6500
$ws.Range("L340").Value = 7000
$ws.Range("M340").Value = 6761
$ws.Range("N340").Value = '$/docena de atados (3 kilos)'
$ws.Range("O340").Value = "Provincia de Quillota"
$ws.Range("P340").Value = 2254
$ws.Range("Q340").Value = 3
$ws.Range("R340").Value = "Hortaliza"
